# Fix Training Data Issue (#48)
# The "Date" column (BF) held the literal folder/file-name string
# "6-22-2011-12" for every data row. Replace it with the correct
# ISO-style date "2012-06-22" (the data was captured one day off
# because of how NBA stats were shown).
#
# NumberFormat is temporarily switched to Text ("@") before writing the
# value so Excel's automatic date-recognition does not turn the literal
# string into a serial date number; the style is then reset back to the
# workbook's default "Normal" style so the cell's formatting is left
# exactly as it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("BF2:BF31")
$rng.NumberFormat = "@"

for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 58).Value = "2012-06-22"
}

$rng.Style = "Normal"
